$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: safe to assign directly.
$ws.Range("D2").Value = "29.496.62"
$ws.Range("E2").Value = "  +3.08%  "
$ws.Range("D3").Value = "1.603.02"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  +8.21%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +9.20%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "1.833.02"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").Value = "1.601.30"
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("D15").Value = "29.523.85"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("E16").Value = "  +4.15%  "
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("E18").Value = "  +3.24%  "
$ws.Range("E19").Value = "  +5.67%  "
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("E21").Value = "  +2.93%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("E27").Value = "  +6.17%  "
$ws.Range("E28").Value = "  +3.61%  "
$ws.Range("E29").Value = "  +2.77%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +2.68%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("E34").Value = "  +3.61%  "
$ws.Range("D35").Value = "1.413.57"
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("E38").Value = "  +5.30%  "
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("E41").Value = "  +3.43%  "
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("E44").Value = "  +21.62%  "
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("E47").Value = "  +2.57%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "1.743.93"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("E51").Value = "  -4.41%  "

# Numeric-looking values that must stay text: round-trip through a
# Text-formatted helper sheet + PasteSpecial so no Number conversion
# and no style/format drift happens on the target cells.
$helperSheet = $wb.Worksheets.Add($null, $ws)
$h = $helperSheet.Cells.Item(1, 1)
$h.NumberFormat = "@"
$h.Value = "212.96"
$h.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(2, 1)
$h.NumberFormat = "@"
$h.Value = "0.526"
$h.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(3, 1)
$h.NumberFormat = "@"
$h.Value = "0.999"
$h.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(4, 1)
$h.NumberFormat = "@"
$h.Value = "26.83"
$h.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(5, 1)
$h.NumberFormat = "@"
$h.Value = "43.48"
$h.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(6, 1)
$h.NumberFormat = "@"
$h.Value = "0.0910"
$h.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(7, 1)
$h.NumberFormat = "@"
$h.Value = "3.71"
$h.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(8, 1)
$h.NumberFormat = "@"
$h.Value = "63.41"
$h.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(9, 1)
$h.NumberFormat = "@"
$h.Value = "242.79"
$h.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(10, 1)
$h.NumberFormat = "@"
$h.Value = "7.61"
$h.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(11, 1)
$h.NumberFormat = "@"
$h.Value = "3.99"
$h.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(12, 1)
$h.NumberFormat = "@"
$h.Value = "154.43"
$h.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(13, 1)
$h.NumberFormat = "@"
$h.Value = "0.109"
$h.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(14, 1)
$h.NumberFormat = "@"
$h.Value = "15.30"
$h.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(15, 1)
$h.NumberFormat = "@"
$h.Value = "2.80"
$h.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(16, 1)
$h.NumberFormat = "@"
$h.Value = "1.96"
$h.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(17, 1)
$h.NumberFormat = "@"
$h.Value = "0.998"
$h.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(18, 1)
$h.NumberFormat = "@"
$h.Value = "52.74"
$h.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(19, 1)
$h.NumberFormat = "@"
$h.Value = "65.63"
$h.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$h = $helperSheet.Cells.Item(20, 1)
$h.NumberFormat = "@"
$h.Value = "86.35"
$h.Copy()
$ws.Range("D50").PasteSpecial(-4163)

$excel.CutCopyMode = $false
$excel.DisplayAlerts = $false
$helperSheet.Delete()
$excel.DisplayAlerts = $true
